$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper note: Price column (D) cells are stored as plain TEXT in this sheet
# even when they look like numbers (e.g. "114.70"). Assigning such a string
# via .Value would make Excel auto-convert it to a numeric cell, which
# doesn't match the source data (still text). Prefixing with a single quote
# forces Excel to keep it as text (quotePrefix); we then reset the cell's
# Style back to "Normal" so no stray number-format style sticks to the cell.

function Set-TextValue($rangeAddr, $text) {
    $r = $ws.Range($rangeAddr)
    $r.Value = "'" + $text
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "50.111.40"
$ws.Range("E2").Value = "  +4.21%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.660.69"
$ws.Range("E3").Value = "  +6.82%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.11%  "

# Row 5 - Solana
Set-TextValue "D5" "114.77"
$ws.Range("E5").Value = "  +8.61%  "

# Row 6 - BNB
Set-TextValue "D6" "326.56"
$ws.Range("E6").Value = "  +2.82%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  +2.19%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.07%  "

# Row 10 - Avalanche
Set-TextValue "D10" "41.47"
$ws.Range("E10").Value = "  +6.41%  "

# Row 11 - Chainlink
Set-TextValue "D11" "20.19"
$ws.Range("E11").Value = "  -0.07%  "

# Row 12 - Dogecoin
$ws.Range("E12").Value = "  +3.13%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +0.51%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  +4.48%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "3.077.74"
$ws.Range("E15").Value = "  +6.74%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "2.656.57"
$ws.Range("E16").Value = "  +6.43%  "

# Row 17 - Polygon
$ws.Range("E17").Value = "  +6.06%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "50.054.33"

# Row 19 - InternetComputer(DFINITY)
Set-TextValue "D19" "13.31"
$ws.Range("E19").Value = "  +4.18%  "

# Row 20 - Uniswap
Set-TextValue "D20" "6.79"
$ws.Range("E20").Value = "  +3.21%  "

# Row 21 - ImmutableX
Set-TextValue "D21" "2.94"
$ws.Range("E21").Value = "  -1.71%  "

# Row 22 - ShibaInu
$ws.Range("E22").Value = "  +3.34%  "

# Row 23 - Litecoin
$ws.Range("E23").Value = "  +2.11%  "

# Row 24 - BitcoinCash
Set-TextValue "D24" "276.71"
$ws.Range("E24").Value = "  +2.51%  "

# Row 25 - PancakeSwap
$ws.Range("E25").Value = "  +3.47%  "

# Row 26 - EthereumClassic
Set-TextValue "D26" "26.99"
$ws.Range("E26").Value = "  +4.95%  "

# Row 27 - Dai
$ws.Range("E27").Value = "  +0.10%  "

# Row 28 - Cosmos
Set-TextValue "D28" "10.07"
$ws.Range("E28").Value = "  +3.48%  "

# Row 29 - InjectiveProtocol
Set-TextValue "D29" "36.89"
$ws.Range("E29").Value = "  +6.67%  "

# Row 30 - Toncoin
$ws.Range("E30").Value = "  +1.47%  "

# Row 31 - Kaspa
$ws.Range("E31").Value = "  +2.54%  "

# Row 32 - OKB
Set-TextValue "D32" "50.24"
$ws.Range("E32").Value = "  +2.17%  "

# Row 33 - Filecoin
Set-TextValue "D33" "5.52"
$ws.Range("E33").Value = "  +4.58%  "

# Row 34 - Celestia
Set-TextValue "D34" "19.77"
$ws.Range("E34").Value = "  +3.51%  "

# Row 35 - Hedera
Set-TextValue "D35" "0.0815"
$ws.Range("E35").Value = "  +5.70%  "

# Row 36 - FirstDigitalUSD
$ws.Range("E36").Value = "  -0.03%  "

# Row 37 - RenderToken
Set-TextValue "D37" "5.03"
$ws.Range("E37").Value = "  +9.88%  "

# Row 38 - ARBITRUM
Set-TextValue "D38" "2.09"
$ws.Range("E38").Value = "  +7.38%  "

# Row 39 - LidoDAOToken
$ws.Range("E39").Value = "  +9.04%  "

# Row 40 - was Stellar, now Monero
$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D40" "124.54"
$ws.Range("E40").Value = "  +1.61%  "

# Row 41 - was Monero, now Stellar
$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D41" "0.113"
$ws.Range("E41").Value = "  +2.59%  "

# Row 42 - EnergySwap
Set-TextValue "D42" "22.29"
$ws.Range("E42").Value = "  +0.41%  "

# Row 43 - WEMIXToken
Set-TextValue "D43" "2.23"
$ws.Range("E43").Value = "  +0.24%  "

# Row 44 - VeChain
$ws.Range("E44").Value = "  +5.05%  "

# Row 45 - Maker
$ws.Range("D45").Value = "2.107.59"
$ws.Range("E45").Value = "  +5.41%  "

# Row 46 - NEARProtocol
Set-TextValue "D46" "3.35"
$ws.Range("E46").Value = "  +5.47%  "

# Row 47 - ApeXProtocol
Set-TextValue "D47" "2.28"
$ws.Range("E47").Value = "  +14.28%  "

# Row 48 - Stacks
$ws.Range("E48").Value = "  +4.81%  "

# Row 49 - FraxShare
$ws.Range("E49").Value = "  +2.36%  "

# Row 51 - MultiversX
Set-TextValue "D51" "60.18"
$ws.Range("E51").Value = "  +6.34%  "
